# Update the auto-updating "datetimeFigureOut" date placeholder text
# (shown in the footer area of the Slide Master, every slide layout, and
# the Notes Master) from "8/25/25" to "9/2/25" -- the date PowerPoint
# re-stamps these fields with whenever the deck is edited/saved.

$p = $ppt.ActivePresentation
$newDate = "9/2/25"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $isDate = $false
            if ($sh.Type -eq 14) {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDate = $true
                }
            }
            if ($isDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master footer date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's footer date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes Master date placeholder.
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
